# Updates the "cryptos" list (Price / Volume(1h) columns) to the latest
# scrape, as produced by the "Updated cryptos list ... with GitHub Actions"
# workflow run. Column D ("Price") cells are plain text in this workbook
# (no number format applied), so any value that Excel's COM layer would
# otherwise auto-convert to a real number is entered with a leading
# apostrophe to force it back to text - this mirrors exactly what typing
# the same value into Excel by hand would require, and the apostrophe
# itself is not stored as part of the cell text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.050.86"
$ws.Range("E2").Value = "  +3.01%  "

$ws.Range("D3").Value = "3.734.57"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'600.74"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").Value = "'167.31"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7").Value = "3.733.74"
$ws.Range("E7").Value = "  +1.33%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("D11").Value = "'6.40"
$ws.Range("E11").Value = "  +4.00%  "

$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "'37.79"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("E14").Value = "  +1.73%  "

$ws.Range("D15").Value = "4.352.16"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").Value = "3.729.20"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").Value = "69.000.06"
$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").Value = "'16.95"
$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("D21").Value = "'10.72"
$ws.Range("E21").Value = "  +16.44%  "

$ws.Range("D22").Value = "'491.51"
$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").Value = "'0.722"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'0.0000150"
$ws.Range("E24").Value = "  +6.27%  "

$ws.Range("D25").Value = "'84.55"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").Value = "'12.27"
$ws.Range("E27").Value = "  +0.66%  "

$ws.Range("D28").Value = "'10.10"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("E31").Value = "  +6.64%  "

$ws.Range("D32").Value = "'8.07"
$ws.Range("E32").Value = "  +5.80%  "

$ws.Range("E33").Value = "  +1.08%  "

$ws.Range("D34").Value = "3.875.17"
$ws.Range("E34").Value = "  +0.89%  "

# Rows 35/36 swap places (Hedera now ranks above RenzoRestakedETH) and get
# freshly scraped price/volume figures.
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.108"
$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.663.59"
$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("D38").Value = "'1.02"
$ws.Range("E38").Value = "  +2.02%  "

$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("D40").Value = "'0.134"
$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("D42").Value = "'2.95"
$ws.Range("E42").Value = "  +5.21%  "

$ws.Range("D43").Value = "'430.71"
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("D44").Value = "'48.61"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("E46").Value = "  +0.45%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").Value = "'40.00"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").Value = "'141.50"
$ws.Range("E49").Value = "  +0.62%  "

$ws.Range("D50").Value = "2.760.83"
$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("E51").Value = "  +0.92%  "
